$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()
$wsOverview.Range("A2").Value = '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md'
$wsOverview.Range("B2").Value = 'Handed back: in sync with en-US'
$wsOverview.Range("C2").Value = 'Handed back: in sync with en-US'
$wsOverview.Range("D2").Value = '2016-03-24 10:20:14'
$rngA2 = $wsOverview.Range("A2")
$wsOverview.Hyperlinks.Add($rngA2, 'https://github.com/OpenLocalizationTest/oltest/blob/2dd91785ba39776e11c286ad0666c6995c7d7834/e2e/569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md', "", "", '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md')

$wsOverview.Range("A3").Value = 'ffff59a32bda-caa8-4d76-8caa-8e980c81ed1f.md'
$wsOverview.Range("B3").Value = 'Handed back: in sync with en-US'
$wsOverview.Range("C3").Value = 'Handed back: in sync with en-US'
$wsOverview.Range("D3").Value = '2016-03-24 10:18:19'
$rngA3 = $wsOverview.Range("A3")
$wsOverview.Hyperlinks.Add($rngA3, 'https://github.com/OpenLocalizationTest/oltest/blob/2dd91785ba39776e11c286ad0666c6995c7d7834/e2e/ffff59a32bda-caa8-4d76-8caa-8e980c81ed1f.md', "", "", 'ffff59a32bda-caa8-4d76-8caa-8e980c81ed1f.md')

$wsOverview.Range("A4").Value = 'ffffff04f0fc7d-ba6c-4eae-90de-5851628c71af.md'
$wsOverview.Range("B4").Value = 'Handed back: in sync with en-US'
$wsOverview.Range("C4").Value = 'Handed back: in sync with en-US'
$wsOverview.Range("D4").Value = '2016-03-24 10:18:19'
$rngA4 = $wsOverview.Range("A4")
$wsOverview.Hyperlinks.Add($rngA4, 'https://github.com/OpenLocalizationTest/oltest/blob/2dd91785ba39776e11c286ad0666c6995c7d7834/e2e/ffffff04f0fc7d-ba6c-4eae-90de-5851628c71af.md', "", "", 'ffffff04f0fc7d-ba6c-4eae-90de-5851628c71af.md')

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Range("A2").Value = '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md'
$wsZhCn.Range("B2").Value = '.md'
$wsZhCn.Range("C2").Value = 'Handed back: in sync with en-US'
$wsZhCn.Range("D2").Value = '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.c1be3fb36d390a9cd3cae817e92c27c732a6249a.zh-cn.xlf'
$wsZhCn.Range("E2").Value = '2016-03-24 10:20:09'
$wsZhCn.Range("F2").Value = '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md'
$wsZhCn.Range("G2").Value = '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.c1be3fb36d390a9cd3cae817e92c27c732a6249a.zh-cn.xlf'
$wsZhCn.Range("H2").Value = '2016-03-24 10:21:14'
$wsZhCn.Range("I2").ClearContents()
$wsZhCn.Range("J2").Value = 'Include'
$wsZhCn.Range("K2").ClearContents()
$wsZhCn.Range("L2").ClearContents()
$rngA2 = $wsZhCn.Range("A2")
$wsZhCn.Hyperlinks.Add($rngA2, 'https://github.com/OpenLocalizationTest/oltest/blob/2dd91785ba39776e11c286ad0666c6995c7d7834/e2e/569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md', "", "", '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md')
$rngD2 = $wsZhCn.Range("D2")
$wsZhCn.Hyperlinks.Add($rngD2, 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7aef5fd22ef3ce2ffc3b75c5f74706af6aed79f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/569b6b7c-1d47-4301-bdc3-1b1af6ef4366.c1be3fb36d390a9cd3cae817e92c27c732a6249a.zh-cn.xlf', "", "", '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.c1be3fb36d390a9cd3cae817e92c27c732a6249a.zh-cn.xlf')
$rngF2 = $wsZhCn.Range("F2")
$wsZhCn.Hyperlinks.Add($rngF2, 'https://github.com/OpenLocalizationTest/oltest/blob/2dd91785ba39776e11c286ad0666c6995c7d7834/e2e/569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md', "", "", '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md')
$rngG2 = $wsZhCn.Range("G2")
$wsZhCn.Hyperlinks.Add($rngG2, 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7aef5fd22ef3ce2ffc3b75c5f74706af6aed79f1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/569b6b7c-1d47-4301-bdc3-1b1af6ef4366.c1be3fb36d390a9cd3cae817e92c27c732a6249a.zh-cn.xlf', "", "", '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.c1be3fb36d390a9cd3cae817e92c27c732a6249a.zh-cn.xlf')

$wsZhCn.Range("A3").Value = 'ffff59a32bda-caa8-4d76-8caa-8e980c81ed1f.md'
$wsZhCn.Range("B3").Value = '.md'
$wsZhCn.Range("C3").Value = 'Handed back: in sync with en-US'
$wsZhCn.Range("D3").Value = '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf'
$wsZhCn.Range("E3").Value = '2016-03-24 10:18:15'
$wsZhCn.Range("F3").Value = '437cfbd5-767f-4178-a01b-f91116985aef.md'
$wsZhCn.Range("G3").Value = '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf'
$wsZhCn.Range("H3").Value = '2016-03-24 10:18:40'
$wsZhCn.Range("I3").ClearContents()
$wsZhCn.Range("J3").Value = 'Include'
$wsZhCn.Range("K3").ClearContents()
$wsZhCn.Range("L3").ClearContents()
$rngA3 = $wsZhCn.Range("A3")
$wsZhCn.Hyperlinks.Add($rngA3, 'https://github.com/OpenLocalizationTest/oltest/blob/2dd91785ba39776e11c286ad0666c6995c7d7834/e2e/ffff59a32bda-caa8-4d76-8caa-8e980c81ed1f.md', "", "", 'ffff59a32bda-caa8-4d76-8caa-8e980c81ed1f.md')
$rngD3 = $wsZhCn.Range("D3")
$wsZhCn.Hyperlinks.Add($rngD3, 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e8b4a8d93646e1dd11c6dd34ae4da2498d871b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf', "", "", '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf')
$rngF3 = $wsZhCn.Range("F3")
$wsZhCn.Hyperlinks.Add($rngF3, 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/842700f175fd9230bedc1e8136cf12b24ad7b962/e2e/437cfbd5-767f-4178-a01b-f91116985aef.md', "", "", '437cfbd5-767f-4178-a01b-f91116985aef.md')
$rngG3 = $wsZhCn.Range("G3")
$wsZhCn.Hyperlinks.Add($rngG3, 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/73474d69fe7de64797aee915a4b365f07fd3bcaa/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf', "", "", '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf')

$wsZhCn.Range("A4").Value = 'ffffff04f0fc7d-ba6c-4eae-90de-5851628c71af.md'
$wsZhCn.Range("B4").Value = '.md'
$wsZhCn.Range("C4").Value = 'Handed back: in sync with en-US'
$wsZhCn.Range("D4").Value = '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf'
$wsZhCn.Range("E4").Value = '2016-03-24 10:18:15'
$wsZhCn.Range("F4").Value = '437cfbd5-767f-4178-a01b-f91116985aef.md'
$wsZhCn.Range("G4").Value = '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf'
$wsZhCn.Range("H4").Value = '2016-03-24 10:18:40'
$wsZhCn.Range("I4").ClearContents()
$wsZhCn.Range("J4").Value = 'Include'
$wsZhCn.Range("K4").ClearContents()
$wsZhCn.Range("L4").ClearContents()
$rngA4 = $wsZhCn.Range("A4")
$wsZhCn.Hyperlinks.Add($rngA4, 'https://github.com/OpenLocalizationTest/oltest/blob/2dd91785ba39776e11c286ad0666c6995c7d7834/e2e/ffffff04f0fc7d-ba6c-4eae-90de-5851628c71af.md', "", "", 'ffffff04f0fc7d-ba6c-4eae-90de-5851628c71af.md')
$rngD4 = $wsZhCn.Range("D4")
$wsZhCn.Hyperlinks.Add($rngD4, 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e8b4a8d93646e1dd11c6dd34ae4da2498d871b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf', "", "", '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf')
$rngF4 = $wsZhCn.Range("F4")
$wsZhCn.Hyperlinks.Add($rngF4, 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/842700f175fd9230bedc1e8136cf12b24ad7b962/e2e/437cfbd5-767f-4178-a01b-f91116985aef.md', "", "", '437cfbd5-767f-4178-a01b-f91116985aef.md')
$rngG4 = $wsZhCn.Range("G4")
$wsZhCn.Hyperlinks.Add($rngG4, 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/73474d69fe7de64797aee915a4b365f07fd3bcaa/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf', "", "", '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf')

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Range("A2").Value = '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md'
$wsDeDe.Range("B2").Value = '.md'
$wsDeDe.Range("C2").Value = 'Handed back: in sync with en-US'
$wsDeDe.Range("D2").Value = '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.c1be3fb36d390a9cd3cae817e92c27c732a6249a.de-de.xlf'
$wsDeDe.Range("E2").Value = '2016-03-24 10:20:14'
$wsDeDe.Range("F2").Value = '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md'
$wsDeDe.Range("G2").Value = '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.c1be3fb36d390a9cd3cae817e92c27c732a6249a.de-de.xlf'
$wsDeDe.Range("H2").Value = '2016-03-24 10:21:22'
$wsDeDe.Range("I2").ClearContents()
$wsDeDe.Range("J2").Value = 'Include'
$wsDeDe.Range("K2").ClearContents()
$wsDeDe.Range("L2").ClearContents()
$rngA2 = $wsDeDe.Range("A2")
$wsDeDe.Hyperlinks.Add($rngA2, 'https://github.com/OpenLocalizationTest/oltest/blob/2dd91785ba39776e11c286ad0666c6995c7d7834/e2e/569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md', "", "", '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md')
$rngD2 = $wsDeDe.Range("D2")
$wsDeDe.Hyperlinks.Add($rngD2, 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e7c3a700e84e0ef111c0f27e70d9010fd4bbae6d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/569b6b7c-1d47-4301-bdc3-1b1af6ef4366.c1be3fb36d390a9cd3cae817e92c27c732a6249a.de-de.xlf', "", "", '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.c1be3fb36d390a9cd3cae817e92c27c732a6249a.de-de.xlf')
$rngF2 = $wsDeDe.Range("F2")
$wsDeDe.Hyperlinks.Add($rngF2, 'https://github.com/OpenLocalizationTest/oltest/blob/2dd91785ba39776e11c286ad0666c6995c7d7834/e2e/569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md', "", "", '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.md')
$rngG2 = $wsDeDe.Range("G2")
$wsDeDe.Hyperlinks.Add($rngG2, 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e7c3a700e84e0ef111c0f27e70d9010fd4bbae6d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/569b6b7c-1d47-4301-bdc3-1b1af6ef4366.c1be3fb36d390a9cd3cae817e92c27c732a6249a.de-de.xlf', "", "", '569b6b7c-1d47-4301-bdc3-1b1af6ef4366.c1be3fb36d390a9cd3cae817e92c27c732a6249a.de-de.xlf')

$wsDeDe.Range("A3").Value = 'ffff59a32bda-caa8-4d76-8caa-8e980c81ed1f.md'
$wsDeDe.Range("B3").Value = '.md'
$wsDeDe.Range("C3").Value = 'Handed back: in sync with en-US'
$wsDeDe.Range("D3").Value = '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf'
$wsDeDe.Range("E3").Value = '2016-03-24 10:18:19'
$wsDeDe.Range("F3").Value = '437cfbd5-767f-4178-a01b-f91116985aef.md'
$wsDeDe.Range("G3").Value = '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf'
$wsDeDe.Range("H3").Value = '2016-03-24 10:18:47'
$wsDeDe.Range("I3").ClearContents()
$wsDeDe.Range("J3").Value = 'Include'
$wsDeDe.Range("K3").ClearContents()
$wsDeDe.Range("L3").ClearContents()
$rngA3 = $wsDeDe.Range("A3")
$wsDeDe.Hyperlinks.Add($rngA3, 'https://github.com/OpenLocalizationTest/oltest/blob/2dd91785ba39776e11c286ad0666c6995c7d7834/e2e/ffff59a32bda-caa8-4d76-8caa-8e980c81ed1f.md', "", "", 'ffff59a32bda-caa8-4d76-8caa-8e980c81ed1f.md')
$rngD3 = $wsDeDe.Range("D3")
$wsDeDe.Hyperlinks.Add($rngD3, 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/23cd3b2a47102d7d34a44624085eb5132c4bf97b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf', "", "", '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf')
$rngF3 = $wsDeDe.Range("F3")
$wsDeDe.Hyperlinks.Add($rngF3, 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a63355ec0552604d08d425ddb3f232501c2dd24c/e2e/437cfbd5-767f-4178-a01b-f91116985aef.md', "", "", '437cfbd5-767f-4178-a01b-f91116985aef.md')
$rngG3 = $wsDeDe.Range("G3")
$wsDeDe.Hyperlinks.Add($rngG3, 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/ef8d2d15ff0848e2ed834cf7dc695c05a3b341d7/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf', "", "", '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf')

$wsDeDe.Range("A4").Value = 'ffffff04f0fc7d-ba6c-4eae-90de-5851628c71af.md'
$wsDeDe.Range("B4").Value = '.md'
$wsDeDe.Range("C4").Value = 'Handed back: in sync with en-US'
$wsDeDe.Range("D4").Value = '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf'
$wsDeDe.Range("E4").Value = '2016-03-24 10:18:19'
$wsDeDe.Range("F4").Value = '437cfbd5-767f-4178-a01b-f91116985aef.md'
$wsDeDe.Range("G4").Value = '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf'
$wsDeDe.Range("H4").Value = '2016-03-24 10:18:47'
$wsDeDe.Range("I4").ClearContents()
$wsDeDe.Range("J4").Value = 'Include'
$wsDeDe.Range("K4").ClearContents()
$wsDeDe.Range("L4").ClearContents()
$rngA4 = $wsDeDe.Range("A4")
$wsDeDe.Hyperlinks.Add($rngA4, 'https://github.com/OpenLocalizationTest/oltest/blob/2dd91785ba39776e11c286ad0666c6995c7d7834/e2e/ffffff04f0fc7d-ba6c-4eae-90de-5851628c71af.md', "", "", 'ffffff04f0fc7d-ba6c-4eae-90de-5851628c71af.md')
$rngD4 = $wsDeDe.Range("D4")
$wsDeDe.Hyperlinks.Add($rngD4, 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/23cd3b2a47102d7d34a44624085eb5132c4bf97b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf', "", "", '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf')
$rngF4 = $wsDeDe.Range("F4")
$wsDeDe.Hyperlinks.Add($rngF4, 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/a63355ec0552604d08d425ddb3f232501c2dd24c/e2e/437cfbd5-767f-4178-a01b-f91116985aef.md', "", "", '437cfbd5-767f-4178-a01b-f91116985aef.md')
$rngG4 = $wsDeDe.Range("G4")
$wsDeDe.Hyperlinks.Add($rngG4, 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/ef8d2d15ff0848e2ed834cf7dc695c05a3b341d7/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf', "", "", '437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf')

